# Update row 2 of Sheet1:
#   A2: "Full Evaluation"   -> "**Overall Score"
#   B2: (long evaluation text) -> "8"  (kept as text, not converted to a number)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "**Overall Score"

# Setting Value = "8" directly would be auto-coerced to the number 8 by
# Excel's smart-typing. Force a text number-format first so the literal
# stays a string, then clear the temporary formatting so no stray style
# is left behind on the cell.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "8"
$ws.Range("B2").ClearFormats()
